$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log_Muestras")

# Apply per-row corrections to percentil_densidad_75 (G), densidad (K),
# pasa_densidad (O), is_filtrada (D) and timestamp (Z) columns, as produced
# by the corrected PCSMOTE density cache / fingerprint logic.

$ws.Range("G2").Value = 0.2857142857142857
$ws.Range("K2").Value = 0.2857142857142857
$ws.Range("O2").Value = $True
$ws.Range("Z2").Value = "2025-10-19T23:55:39.462558"

$ws.Range("G3").Value = 0.2857142857142857
$ws.Range("K3").Value = 0.2857142857142857
$ws.Range("O3").Value = $True
$ws.Range("Z3").Value = "2025-10-19T23:55:39.462558"

$ws.Range("G4").Value = 0.2857142857142857
$ws.Range("K4").Value = 0.2857142857142857
$ws.Range("Z4").Value = "2025-10-19T23:55:39.462558"

$ws.Range("G5").Value = 0.2857142857142857
$ws.Range("K5").Value = 0.4285714285714285
$ws.Range("Z5").Value = "2025-10-19T23:55:39.462558"

$ws.Range("D6").Value = $True
$ws.Range("G6").Value = 0.2857142857142857
$ws.Range("K6").Value = 0.2857142857142857
$ws.Range("O6").Value = $True
$ws.Range("Z6").Value = "2025-10-19T23:55:39.463562"

$ws.Range("G7").Value = 0.2857142857142857
$ws.Range("K7").Value = 0.2857142857142857
$ws.Range("O7").Value = $True
$ws.Range("Z7").Value = "2025-10-19T23:55:39.463562"

$ws.Range("G8").Value = 0.2857142857142857
$ws.Range("K8").Value = 0.2857142857142857
$ws.Range("O8").Value = $True
$ws.Range("Z8").Value = "2025-10-19T23:55:39.463562"

$ws.Range("G9").Value = 0.2857142857142857
$ws.Range("K9").Value = 0.2857142857142857
$ws.Range("O9").Value = $True
$ws.Range("Z9").Value = "2025-10-19T23:55:39.463562"

$ws.Range("G10").Value = 0.2857142857142857
$ws.Range("K10").Value = 0.2857142857142857
$ws.Range("O10").Value = $True
$ws.Range("Z10").Value = "2025-10-19T23:55:39.464561"

$ws.Range("G11").Value = 0.2857142857142857
$ws.Range("K11").Value = 0.2857142857142857
$ws.Range("O11").Value = $True
$ws.Range("Z11").Value = "2025-10-19T23:55:39.464561"

$ws.Range("G12").Value = 0.2857142857142857
$ws.Range("K12").Value = 0.2857142857142857
$ws.Range("O12").Value = $True
$ws.Range("Z12").Value = "2025-10-19T23:55:39.464561"

$ws.Range("G13").Value = 0.2857142857142857
$ws.Range("K13").Value = 0.2857142857142857
$ws.Range("O13").Value = $True
$ws.Range("Z13").Value = "2025-10-19T23:55:39.464561"

$ws.Range("G14").Value = 0.2857142857142857
$ws.Range("K14").Value = 0.4285714285714285
$ws.Range("Z14").Value = "2025-10-19T23:55:39.465556"

$ws.Range("G15").Value = 0.2857142857142857
$ws.Range("K15").Value = 0.5714285714285714
$ws.Range("Z15").Value = "2025-10-19T23:55:39.465556"

$ws.Range("G16").Value = 0.2857142857142857
$ws.Range("K16").Value = 0.2857142857142857
$ws.Range("Z16").Value = "2025-10-19T23:55:39.507562"

$ws.Range("G17").Value = 0.2857142857142857
$ws.Range("K17").Value = 0.2857142857142857
$ws.Range("Z17").Value = "2025-10-19T23:55:39.508557"

$ws.Range("G18").Value = 0.2857142857142857
$ws.Range("K18").Value = 0.2857142857142857
$ws.Range("Z18").Value = "2025-10-19T23:55:39.508557"

$ws.Range("G19").Value = 0.2857142857142857
$ws.Range("K19").Value = 0.2857142857142857
$ws.Range("Z19").Value = "2025-10-19T23:55:39.509561"

$ws.Range("G20").Value = 0.2857142857142857
$ws.Range("K20").Value = 0.2857142857142857
$ws.Range("Z20").Value = "2025-10-19T23:55:39.509561"

$ws.Range("G21").Value = 0.2857142857142857
$ws.Range("K21").Value = 0.2857142857142857
$ws.Range("Z21").Value = "2025-10-19T23:55:39.509561"

$ws.Range("G22").Value = 0.2857142857142857
$ws.Range("K22").Value = 0.2857142857142857
$ws.Range("Z22").Value = "2025-10-19T23:55:39.510560"

$ws.Range("G23").Value = 0.2857142857142857
$ws.Range("K23").Value = 0.2857142857142857
$ws.Range("Z23").Value = "2025-10-19T23:55:39.510560"

$ws.Range("G24").Value = 0.2857142857142857
$ws.Range("K24").Value = 0.2857142857142857
$ws.Range("Z24").Value = "2025-10-19T23:55:39.510560"

$ws.Range("G25").Value = 0.2857142857142857
$ws.Range("K25").Value = 0.2857142857142857
$ws.Range("Z25").Value = "2025-10-19T23:55:39.510560"

$ws.Range("G26").Value = 0.2857142857142857
$ws.Range("K26").Value = 0.7142857142857143
$ws.Range("Z26").Value = "2025-10-19T23:55:39.549563"

$ws.Range("G27").Value = 0.2857142857142857
$ws.Range("K27").Value = 0.7142857142857143
$ws.Range("Z27").Value = "2025-10-19T23:55:39.549563"

$ws.Range("G28").Value = 0.2857142857142857
$ws.Range("K28").Value = 0.2857142857142857
$ws.Range("O28").Value = $True
$ws.Range("Z28").Value = "2025-10-19T23:55:39.549563"

$ws.Range("G29").Value = 0.2857142857142857
$ws.Range("K29").Value = 0.2857142857142857
$ws.Range("O29").Value = $True
$ws.Range("Z29").Value = "2025-10-19T23:55:39.550561"

$ws.Range("G30").Value = 0.2857142857142857
$ws.Range("K30").Value = 0.2857142857142857
$ws.Range("O30").Value = $True
$ws.Range("Z30").Value = "2025-10-19T23:55:39.550561"

$ws.Range("G31").Value = 0.2857142857142857
$ws.Range("K31").Value = 0.2857142857142857
$ws.Range("O31").Value = $True
$ws.Range("Z31").Value = "2025-10-19T23:55:39.550561"

$ws.Range("G32").Value = 0.2857142857142857
$ws.Range("K32").Value = 0.2857142857142857
$ws.Range("O32").Value = $True
$ws.Range("Z32").Value = "2025-10-19T23:55:39.550561"

$ws.Range("G33").Value = 0.2857142857142857
$ws.Range("K33").Value = 0.2857142857142857
$ws.Range("O33").Value = $True
$ws.Range("Z33").Value = "2025-10-19T23:55:39.551562"

$ws.Range("G34").Value = 0.2857142857142857
$ws.Range("K34").Value = 0.2857142857142857
$ws.Range("O34").Value = $True
$ws.Range("Z34").Value = "2025-10-19T23:55:39.551562"

$ws.Range("G35").Value = 0.2857142857142857
$ws.Range("K35").Value = 0.2857142857142857
$ws.Range("O35").Value = $True
$ws.Range("Z35").Value = "2025-10-19T23:55:39.551562"

$ws.Range("G36").Value = 0.2857142857142857
$ws.Range("K36").Value = 0.4285714285714285
$ws.Range("Z36").Value = "2025-10-19T23:55:39.551562"

$ws.Range("G37").Value = 0.2857142857142857
$ws.Range("K37").Value = 0.2857142857142857
$ws.Range("O37").Value = $True
$ws.Range("Z37").Value = "2025-10-19T23:55:39.551562"

$ws.Range("G38").Value = 0.2857142857142857
$ws.Range("K38").Value = 0.4285714285714285
$ws.Range("Z38").Value = "2025-10-19T23:55:39.557554"

$ws.Range("G39").Value = 0.2857142857142857
$ws.Range("K39").Value = 0.2857142857142857
$ws.Range("Z39").Value = "2025-10-19T23:55:39.557554"

$ws.Range("G40").Value = 0.2857142857142857
$ws.Range("K40").Value = 0.2857142857142857
$ws.Range("O40").Value = $True
$ws.Range("Z40").Value = "2025-10-19T23:55:39.557554"

$ws.Range("G41").Value = 0.2857142857142857
$ws.Range("K41").Value = 0.5714285714285714
$ws.Range("Z41").Value = "2025-10-19T23:55:39.558555"

$ws.Range("G42").Value = 0.2857142857142857
$ws.Range("K42").Value = 0.2857142857142857
$ws.Range("O42").Value = $True
$ws.Range("Z42").Value = "2025-10-19T23:55:39.558555"

$ws.Range("G43").Value = 0.2857142857142857
$ws.Range("K43").Value = 0.2857142857142857
$ws.Range("O43").Value = $True
$ws.Range("Z43").Value = "2025-10-19T23:55:39.558555"

$ws.Range("G44").Value = 0.2857142857142857
$ws.Range("K44").Value = 0.2857142857142857
$ws.Range("O44").Value = $True
$ws.Range("Z44").Value = "2025-10-19T23:55:39.558555"

$ws.Range("G45").Value = 0.2857142857142857
$ws.Range("K45").Value = 0.2857142857142857
$ws.Range("O45").Value = $True
$ws.Range("Z45").Value = "2025-10-19T23:55:39.561559"

$ws.Range("G46").Value = 0.2857142857142857
$ws.Range("K46").Value = 0.2857142857142857
$ws.Range("O46").Value = $True
$ws.Range("Z46").Value = "2025-10-19T23:55:39.562560"

$ws.Range("G47").Value = 0.2857142857142857
$ws.Range("K47").Value = 0.2857142857142857
$ws.Range("O47").Value = $True
$ws.Range("Z47").Value = "2025-10-19T23:55:39.562560"

$ws.Range("G48").Value = 0.2857142857142857
$ws.Range("K48").Value = 0.2857142857142857
$ws.Range("O48").Value = $True
$ws.Range("Z48").Value = "2025-10-19T23:55:39.562560"
